{"js": "// Remove the leading \"NUM \" run from the title paragraph (the paragraph\n// that also holds the \"_GoBack\" bookmark and the \"Telecommunications\n// Contractor\" text), leaving the rest of the paragraph content (including\n// the bookmark) untouched.\nconst body = context.document.body;\nconst results = body.search(\"NUM \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const r = results.items[i];\n  if (r.text === \"NUM \") {\n    r.delete();\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"NUM \"\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nwhile ($rng.Find.Execute()) {\n    $rng.Delete()\n    $rng.Collapse(0)\n}\n"}
